$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B-F (1..5 relative) and I-N (skipping G,H which are unchanged) for rows 2..25.
# Data below corresponds to the updated vm_pu results after switching the
# case to the 380 kV slack-bus voltage setpoint of 1.02 (previously 1.05).

$bfData = @(
    @(1.02,1.030592923649507,1.034187308579446,1.04042822946007,1.051818633936164),
    @(1.02,1.031572355098161,1.034906830788792,1.041324154351524,1.052858884426698),
    @(1.02,1.032206530588527,1.035372523355762,1.041904629852408,1.053532907492634),
    @(1.02,1.032473237321001,1.035568326462688,1.042148840755045,1.053816483982813),
    @(1.02,1.032518024361654,1.035601204151658,1.042189855305786,1.053864110430432),
    @(1.02,1.032210093947377,1.035375139586048,1.04190789230843,1.053536695803411),
    @(1.02,1.030923840165796,1.034430450214997,1.040730855080967,1.052170002134025),
    @(1.02,1.028660526737434,1.032766711428142,1.038662573265564,1.04976873821572),
    @(1.02,1.027153861304967,1.031658246728944,1.037287686446848,1.048172677199668),
    @(1.02,1.026501990295795,1.031178449093043,1.036693299210094,1.047482711136482),
    @(1.02,1.026259935797699,1.031000258346884,1.036472660720264,1.047226599025984),
    @(1.02,1.026311853706698,1.031038479605268,1.036519981920401,1.04728152814494),
    @(1.02,1.026481980372236,1.03116371922339,1.036675058232725,1.047461537315438),
    @(1.02,1.026586811565657,1.03124088715947,1.036770624897913,1.047572469725925),
    @(1.02,1.027197134944961,1.031690093103679,1.037327154024697,1.048218492046836),
    @(1.02,1.027580115626789,1.031971915784173,1.037676504504907,1.048624030328982),
    @(1.02,1.027803552404108,1.032136314981071,1.037880366114427,1.048860683773477),
    @(1.02,1.027879747141774,1.032192373671544,1.037949893125257,1.048941395083362),
    @(1.02,1.027539020183807,1.031941677131893,1.037639013045004,1.04858050854294),
    @(1.02,1.026431880112855,1.031126838509876,1.036629388157523,1.047408524345873),
    @(1.02,1.025736236971603,1.030614676580337,1.03599542705148,1.046672646988206),
    @(1.02,1.026104966802823,1.030886167800603,1.036331422809315,1.047062654823466),
    @(1.02,1.02755758928911,1.031955340625459,1.037655953539635,1.048600173825371),
    @(1.02,1.029245260782329,1.033196710704159,1.039196578943375,1.050388683996512)
)

$inData = @(
    @(1.033150075563958,1.035733313412631,1.036987503238264,1.04321060599651,1.054569081445856,1.037204173517949),
    @(1.033312061061129,1.036354111283547,1.037516400136943,1.043916727251089,1.055421469447415,1.037825852993048),
    @(1.033415605183089,1.036755635272152,1.037858075914158,1.044373753486854,1.055973327151597,1.038227947191772),
    @(1.033458830298964,1.036924393567986,1.038001582442069,1.044565914840021,1.056205400373691,1.038396945143743),
    @(1.033466070096783,1.036952726319668,1.038025669937721,1.044598181165759,1.056244370698122,1.038425318131182),
    @(1.03341618395758,1.036757890396314,1.03785999398293,1.04437632104919,1.055976427842651,1.038230205518469),
    @(1.033205082271905,1.035943150457198,1.037166360985875,1.04344921791211,1.054857086108311,1.037414308555188),
    @(1.032823378242676,1.034506178877643,1.035939872809836,1.041816489780597,1.052887045230777,1.035975296311166),
    @(1.03256240433022,1.033547372811148,1.035119429938369,1.040728690569381,1.051575340653215,1.035015128630081),
    @(1.032447862677774,1.033132012917752,1.034763519457806,1.040257834491298,1.051007761911146,1.034599178877993),
    @(1.032405086135189,1.032977701845578,1.034631221126697,1.040082963623398,1.050796998526994,1.034444648666398),
    @(1.032414272280329,1.033010803346088,1.034659603959118,1.040120472827066,1.050842205231139,1.034477797174835),
    @(1.032444331458596,1.033119258081697,1.034752585617281,1.040243379074989,1.050990338902333,1.034586405928608),
    @(1.032462821356839,1.033186076946247,1.03480986181279,1.040319109145259,1.051081617005636,1.03465331968361),
    @(1.032569973710321,1.033574934925659,1.035143036850413,1.040759943346044,1.051613017441807,1.035042729885956),
    @(1.0326367760258,1.033818804752453,1.035351854419926,1.041036512564854,1.05194645804137,1.035286946035876),
    @(1.032675592240904,1.033961031464316,1.035473591086034,1.04119784688346,1.052140986646228,1.035429374725984),
    @(1.032688802375246,1.034009523949525,1.035515089417254,1.041252860436264,1.052207322351865,1.035477936076083),
    @(1.032629624120849,1.03379264172934,1.035329456790615,1.041006837626633,1.051910679048439,1.035260745858269),
    @(1.032435486144058,1.033087321630053,1.034725207516705,1.040207185523105,1.050946715547816,1.034554424123539),
    @(1.032312089503407,1.032643697742371,1.034344729290487,1.039704563579819,1.050340984650847,1.034110170239053),
    @(1.032377630738728,1.032878886195466,1.034546481063113,1.039970998400937,1.050662060488497,1.034345692686727),
    @(1.032632856217546,1.033804463730296,1.035339577510195,1.041020246403943,1.051926845916813,1.035272584647823),
    @(1.032923206400249,1.03487781927209,1.036257443924391,1.042238471865717,1.053396060121534,1.036347464477602)
)

$bfCols = @(2,3,4,5,6)
$inCols = @(9,10,11,12,13,14)

for ($i = 0; $i -lt $bfData.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $bfCols.Length; $j++) {
        $ws.Cells.Item($row, $bfCols[$j]).Value = $bfData[$i][$j]
    }
    for ($j = 0; $j -lt $inCols.Length; $j++) {
        $ws.Cells.Item($row, $inCols[$j]).Value = $inData[$i][$j]
    }
}

Write-Output "Updated vm_pu values for rows 2-25"
